$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 ("*Doc Ref. No."): 17586324521 -> 13207876678
# D2 ("*Buyer"):        Cloris629k  -> Terencet18o
# E2 ("*Supplier"):     Alicerlzq   -> Julie3396
#
# B2's new value looks numeric, so a plain .Value assignment would make
# Excel store it as a number (losing the shared-string/text type and the
# original "General" style). To keep it as text without touching the
# cell's formatting, write it as a formula that evaluates to a string,
# then convert that formula to its literal value in place.
$ws.Range("B2").Formula = "=""13207876678"""
$ws.Range("B2").Copy()
$ws.Range("B2").PasteSpecial(-4163)

$ws.Range("D2").Value = "Terencet18o"
$ws.Range("E2").Value = "Julie3396"
